$d = $word.ActiveDocument

# 1) Shrink the first run's text from "Apr 2015 " to "Apr " and delete the
#    existing "_GoBack" bookmark so we can re-create it in its new spot.
$d.Content.Find.Execute("Apr 2015 –", $true, $false, $false, $false, $false,
                          $true, 1, $false, "Apr –", 2) | Out-Null

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Re-insert the "_GoBack" bookmark right after "Apr " (before the en dash).
$found = $d.Content.Find.Execute("Apr ")
$r = $d.Content.Find.Parent
$gobackRange = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $gobackRange) | Out-Null
